$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D: "password" ---------------------------------------------

# Header cell D1: reuse the same look as the other header cells (A1:C1)
# by copying C1's format onto D1, then overwrite the text.
$headerSrc = $ws.Range("C1")
$headerDst = $ws.Range("D1")
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)
$headerDst.Value = "password"

# Data cells D2:D27: every teacher gets the same placeholder password.
# Build the look (thin border on all sides + centered text) on D2 first
# so only one new style gets created, then stamp that format onto the
# rest of the column via copy/paste-special (reuses the same style
# instead of minting a new one per cell).
$firstData = $ws.Range("D2")
$firstData.Value = 111111
$firstData.Borders.ColorIndex = 1
$firstData.Borders.LineStyle = 1
$firstData.Borders.Weight = 2
$firstData.HorizontalAlignment = -4108
$firstData.VerticalAlignment = -4108

$restData = $ws.Range("D3:D27")
$restData.Value = 111111
$firstData.Copy()
$restData.PasteSpecial(-4122)

# Column width to roughly match the authored 26.625-character width
# (this COM layer only honours 1/7-character increments).
$ws.Columns.Item(4).ColumnWidth = 25.86

# --- Selection / dimension bookkeeping -------------------------------------
[void]$ws.Range("C9").Select()
